# Implemented getting number of lines for methods and classes.
#
# 1) Update the existing "classFields" sheet: some field rows were
#    re-ordered/re-derived by the analyzer, so their Field Name / Field
#    Modifier / Field Type values move around a bit.
# 2) Add two new sheets "classNumberOfLines" and "methodNumberOfLines"
#    with the computed line-count metrics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) classFields sheet updates
# ---------------------------------------------------------------------
$fields = $wb.Worksheets.Item("classFields")

# Row 2 / Row 3 (pl.piomin.order.OrderApp) swap content
$fields.Range("B2").Value = "LOG"
$fields.Range("C2").Value = "private"
$fields.Range("D2").Value = "org.slf4j.Logger"

$fields.Range("B3").Value = "orderManageService"
$fields.Range("C3").ClearContents()
$fields.Range("D3").Value = "pl.piomin.order.service.OrderManageService"

# Row 10 / Row 11 (pl.piomin.order.OrderControllerTests) swap content
$fields.Range("B10").Value = "factory"
$fields.Range("C10").Value = "private"
$fields.Range("D10").Value = "org.springframework.kafka.core.ConsumerFactory"

$fields.Range("B11").Value = "mapper"
$fields.Range("C11").ClearContents()
$fields.Range("D11").Value = "com.fasterxml.jackson.databind.ObjectMapper"

# Rows 12/13/15 (pl.piomin.order.controller.OrderController) rotate content
$fields.Range("B12").Value = "template"
$fields.Range("D12").Value = "org.springframework.kafka.core.KafkaTemplate"

$fields.Range("B13").Value = "orderGeneratorService"
$fields.Range("D13").Value = "pl.piomin.order.service.OrderGeneratorService"

$fields.Range("B15").Value = "LOG"
$fields.Range("D15").Value = "org.slf4j.Logger"

# ---------------------------------------------------------------------
# 2) New sheet: classNumberOfLines
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$classLines = $wb.Worksheets.Add($null, $lastSheet)
$classLines.Name = "classNumberOfLines"

$classLinesData = @(
    @("Class Name", "Number of Lines"),
    @("pl.piomin.order.OrderAppTest", "6"),
    @("pl.piomin.order.OrderApp", "44"),
    @("pl.piomin.order.KafkaContainerDevMode", "6"),
    @("pl.piomin.order.service.OrderGeneratorService", "24"),
    @("pl.piomin.order.OrderControllerTests", "17"),
    @("pl.piomin.order.controller.OrderController", "34"),
    @("pl.piomin.order.service.OrderManageService", "18")
)

$usedRange = $classLines.Range("A1:B8")
$usedRange.NumberFormat = "@"
for ($i = 0; $i -lt $classLinesData.Length; $i++) {
    $r = $i + 1
    $row = $classLinesData[$i]
    $classLines.Cells.Item($r, 1).Value = $row[0]
    $classLines.Cells.Item($r, 2).Value = $row[1]
}
$usedRange.ClearFormats()

# ---------------------------------------------------------------------
# 3) New sheet: methodNumberOfLines
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$methodLines = $wb.Worksheets.Add($null, $lastSheet2)
$methodLines.Name = "methodNumberOfLines"

$methodLinesData = @(
    @("Class Name", "Method Signature", "Number of Lines"),
    @("pl.piomin.order.OrderAppTest", "main(java.lang.String[])", "3"),
    @("pl.piomin.order.OrderApp", "main(java.lang.String[])", "3"),
    @("pl.piomin.order.OrderApp", "orders()", "3"),
    @("pl.piomin.order.OrderApp", "paymentTopic()", "3"),
    @("pl.piomin.order.OrderApp", "stockTopic()", "3"),
    @("pl.piomin.order.OrderApp", "stream(org.apache.kafka.streams.StreamsBuilder)", "8"),
    @("pl.piomin.order.OrderApp", "table(org.apache.kafka.streams.StreamsBuilder)", "6"),
    @("pl.piomin.order.OrderApp", "taskExecutor()", "8"),
    @("pl.piomin.order.OrderApp", "lambda`$stream`$0(java.lang.Long, pl.piomin.base.domain.Order)", "8"),
    @("pl.piomin.order.KafkaContainerDevMode", "kafka()", "3"),
    @("pl.piomin.order.service.OrderGeneratorService", "generate()", "8"),
    @("pl.piomin.order.OrderControllerTests", "add()", "10"),
    @("pl.piomin.order.controller.OrderController", "create(pl.piomin.base.domain.Order)", "4"),
    @("pl.piomin.order.controller.OrderController", "create()", "4"),
    @("pl.piomin.order.controller.OrderController", "all()", "7"),
    @("pl.piomin.order.controller.OrderController", "lambda`$all`$0(java.util.List, org.apache.kafka.streams.KeyValue)", "7"),
    @("pl.piomin.order.service.OrderManageService", "confirm(pl.piomin.base.domain.Order, pl.piomin.base.domain.Order)", "5")
)

$usedRange2 = $methodLines.Range("A1:C18")
$usedRange2.NumberFormat = "@"
for ($i = 0; $i -lt $methodLinesData.Length; $i++) {
    $r = $i + 1
    $row = $methodLinesData[$i]
    $methodLines.Cells.Item($r, 1).Value = $row[0]
    $methodLines.Cells.Item($r, 2).Value = $row[1]
    $methodLines.Cells.Item($r, 3).Value = $row[2]
}
$usedRange2.ClearFormats()
